$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2 entirely (cell removed from sheet in target)
$ws.Range("C2").ClearContents()

# Update values per diff (small precision corrections from naive forecaster bug fix)
$ws.Range("E2").Value = 4.566338461218034

$ws.Range("C3").Value = 8.604123301398015
$ws.Range("E3").Value = 8.260999835306748

$ws.Range("E4").Value = 4.862860110364853

$ws.Range("C5").Value = 6.334380382529448

$ws.Range("C6").Value = 3.889938592324382
$ws.Range("E6").Value = 3.430035192100678

$ws.Range("C7").Value = 2.513767348245066

$ws.Range("E8").Value = 2.551173534479356

$ws.Range("C10").Value = 2.321003614014905

$ws.Range("E13").Value = 0.9311475558545279

$ws.Range("E14").Value = 2.152035263856322

$ws.Range("E15").Value = -0.1151106897825049

$ws.Range("C18").Value = -0.1883299148263795

$ws.Range("C19").Value = 3.083905204716264
